$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '64.801.37'
$ws.Range('E2').Value = '  -0.29%  '

$ws.Range('D3').Value = '3.144.68'
$ws.Range('E3').Value = '  +0.08%  '

$ws.Range('E4').Value = '  +0.00%  '

Set-TextValue 'D5' '576.55'
$ws.Range('E5').Value = '  +0.99%  '

Set-TextValue 'D6' '148.85'
$ws.Range('E6').Value = '  -0.63%  '

Set-TextValue 'D7' '1.00'
$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('D8').Value = '3.143.22'
$ws.Range('E8').Value = '  +0.03%  '

Set-TextValue 'D9' '0.523'
$ws.Range('E9').Value = '  -0.67%  '

$ws.Range('E10').Value = '  -2.62%  '

Set-TextValue 'D11' '6.07'
$ws.Range('E11').Value = '  -0.98%  '

$ws.Range('E12').Value = '  -1.45%  '

Set-TextValue 'D13' '0.0000258'
$ws.Range('E13').Value = '  +1.83%  '

Set-TextValue 'D14' '36.91'
$ws.Range('E14').Value = '  -1.16%  '

$ws.Range('D15').Value = '3.660.93'
$ws.Range('E15').Value = '  +0.10%  '

$ws.Range('D16').Value = '64.915.30'
$ws.Range('E16').Value = '  -0.19%  '

$ws.Range('D17').Value = '3.145.48'
$ws.Range('E17').Value = '  -0.03%  '

Set-TextValue 'D18' '7.06'
$ws.Range('E18').Value = '  -1.61%  '

$ws.Range('E19').Value = '  +0.07%  '

Set-TextValue 'D20' '503.15'
$ws.Range('E20').Value = '  -1.09%  '

Set-TextValue 'D21' '14.74'
$ws.Range('E21').Value = '  -0.76%  '

$ws.Range('B22').Value = 'InternetComputer(DFINITY)'
$ws.Range('C22').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D22' '15.19'
$ws.Range('E22').Value = '  -2.05%  '

$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D23' '0.709'
$ws.Range('E23').Value = '  -2.77%  '

Set-TextValue 'D24' '7.67'
$ws.Range('E24').Value = '  -1.86%  '

Set-TextValue 'D25' '83.65'
$ws.Range('E25').Value = '  -2.13%  '

Set-TextValue 'D26' '0.996'
$ws.Range('E26').Value = '  -0.23%  '

$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D27' '2.88'
$ws.Range('E27').Value = '  -1.51%  '

$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D28' '8.82'
$ws.Range('E28').Value = '  +1.07%  '

Set-TextValue 'D29' '2.17'
$ws.Range('E29').Value = '  -0.57%  '

Set-TextValue 'D30' '2.79'
$ws.Range('E30').Value = '  +5.51%  '

Set-TextValue 'D31' '27.44'
$ws.Range('E31').Value = '  -1.81%  '

Set-TextValue 'D32' '1.00'
$ws.Range('E32').Value = '  -0.08%  '

$ws.Range('E33').Value = '  -0.37%  '

Set-TextValue 'D34' '6.15'
$ws.Range('E34').Value = '  +2.18%  '

Set-TextValue 'D35' '6.43'
$ws.Range('E35').Value = '  -2.50%  '

Set-TextValue 'D36' '54.56'
$ws.Range('E36').Value = '  -1.93%  '

Set-TextValue 'D37' '0.0888'
$ws.Range('E37').Value = '  +3.82%  '

Set-TextValue 'D38' '477.99'
$ws.Range('E38').Value = '  +1.33%  '

$ws.Range('E39').Value = '  -2.15%  '

Set-TextValue 'D40' '2.95'
$ws.Range('E40').Value = '  -2.40%  '

Set-TextValue 'D41' '8.60'
$ws.Range('E41').Value = '  -0.12%  '

$ws.Range('D42').Value = '3.002.53'
$ws.Range('E42').Value = '  -3.90%  '

Set-TextValue 'D43' '0.115'
$ws.Range('E43').Value = '  -3.05%  '

$ws.Range('E44').Value = '  -3.33%  '

$ws.Range('E45').Value = '  -0.71%  '

Set-TextValue 'D46' '27.99'
$ws.Range('E46').Value = '  -4.01%  '

$ws.Range('E47').Value = '  -0.13%  '

$ws.Range('E48').Value = '  -0.11%  '

$ws.Range('E49').Value = '  -2.14%  '

Set-TextValue 'D50' '2.22'
$ws.Range('E50').Value = '  -4.01%  '

$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D51' '116.84'
$ws.Range('E51').Value = '  -1.25%  '
